# Add a new "pixel_size_mm" column to the im_seq_par parameter sheet.
# The DVF errors in this sheet are supposed isotropic, so they need to be
# multiplied by the pixel size (mm) to get correct values; this commit
# records that pixel size next to the other parameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell (I1): same bold header look as the rest of row 1.
$ws.Range("I1").Value = "pixel_size_mm"
$ws.Range("I1").Font.Bold = $true

# Data cell (I2): the (isotropic) pixel size, in millimetres.
$ws.Range("I2").Value = 1.818

# Reflect where the user ended up after entering the data.
$ws.Range("J6").Select() | Out-Null
